$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 4 with the new certificate entry.
# Order matters for shared-string allocation: A4, C4, E4, B4 (so new
# strings land at sharedStrings indices 13,14,15,16 respectively),
# then the numeric cells D4 and F4.
$ws.Range("A4").Value = "компания легенд"
$ws.Range("C4").Value = "00-00000"
$ws.Range("E4").Value = "РФ"
$ws.Range("B4").Value = "не действует"
$ws.Range("D4").Value = 0
$ws.Range("F4").Value = [DateTime]"2003-10-27"

# Move the active selection from C9 to C8.
$ws.Range("C8").Select() | Out-Null
